$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '35.246.56'
Set-TextValue "E2" '  -0.28%  '
Set-TextValue "D3" '1.905.54'
Set-TextValue "E3" '  -0.34%  '
Set-TextValue "E4" '  +0.23%  '
Set-TextValue "D5" '0.722'
Set-TextValue "E5" '  +9.85%  '
Set-TextValue "D6" '252.88'
Set-TextValue "E6" '  +3.45%  '
Set-TextValue "E7" '  +0.21%  '
Set-TextValue "D8" '40.49'
Set-TextValue "E8" '  -1.89%  '
Set-TextValue "D9" '0.362'
Set-TextValue "E9" '  +3.24%  '
Set-TextValue "D10" '52.26'
Set-TextValue "E10" '  -1.10%  '
Set-TextValue "D11" '0.0762'
Set-TextValue "E11" '  +6.59%  '
Set-TextValue "D12" '0.0988'
Set-TextValue "E12" '  -0.65%  '
Set-TextValue "D13" '2.185.18'
Set-TextValue "E13" '  -0.28%  '
Set-TextValue "D14" '12.72'
Set-TextValue "E14" '  +5.49%  '
Set-TextValue "D15" '0.716'
Set-TextValue "E15" '  +2.16%  '
Set-TextValue "D16" '1.914.48'
Set-TextValue "E16" '  -0.28%  '
Set-TextValue "D17" '4.89'
Set-TextValue "E17" '  +0.56%  '
Set-TextValue "D18" '35.268.99'
Set-TextValue "E18" '  -0.15%  '
Set-TextValue "D19" '74.14'
Set-TextValue "E19" '  +2.99%  '
Set-TextValue "D20" '0.0₃0846'
Set-TextValue "E20" '  +3.13%  '
Set-TextValue "D21" '242.93'
Set-TextValue "E21" '  +1.48%  '
Set-TextValue "D22" '12.96'
Set-TextValue "E22" '  +3.84%  '
Set-TextValue "E23" '  +5.40%  '
Set-TextValue "E24" '  +0.32%  '
Set-TextValue "D25" '2.37'
Set-TextValue "E25" '  +3.65%  '
Set-TextValue "D26" '2.43'
Set-TextValue "E26" '  +3.46%  '
Set-TextValue "D27" '167.12'
Set-TextValue "E27" '  -1.72%  '
Set-TextValue "D28" '8.59'
Set-TextValue "E28" '  +1.68%  '
Set-TextValue "D29" '18.67'
Set-TextValue "E29" '  +1.55%  '
Set-TextValue "E30" '  +4.22%  '
Set-TextValue "D31" '4.126.59'
Set-TextValue "E31" '  +19.39%  '
Set-TextValue "E32" '  +4.38%  '
Set-TextValue "E33" '  +13.76%  '
Set-TextValue "B34" 'TrustWalletToken'
Set-TextValue "C34" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D34" '1.64'
Set-TextValue "E34" '  +23.67%  '
Set-TextValue "B35" 'Hedera'
Set-TextValue "C35" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D35" '0.0581'
Set-TextValue "E35" '  +2.48%  '
Set-TextValue "D36" '4.19'
Set-TextValue "E36" '  +2.01%  '
Set-TextValue "E37" '  +0.10%  '
Set-TextValue "D38" '0.918'
Set-TextValue "E38" '  -1.61%  '
Set-TextValue "D39" '2.02'
Set-TextValue "E39" '  -0.48%  '
Set-TextValue "D40" '0.0217'
Set-TextValue "E40" '  +4.25%  '
Set-TextValue "D41" '17.08'
Set-TextValue "E41" '  +4.38%  '
Set-TextValue "D42" '96.33'
Set-TextValue "E42" '  +6.74%  '
Set-TextValue "E43" '  -0.12%  '
Set-TextValue "D44" '0.0645'
Set-TextValue "E44" '  -3.19%  '
Set-TextValue "D45" '1.336.33'
Set-TextValue "E45" '  -0.27%  '
Set-TextValue "D46" '2.42'
Set-TextValue "E46" '  +2.31%  '
Set-TextValue "E47" '  +0.49%  '
Set-TextValue "D48" '6.72'
Set-TextValue "E48" '  +2.22%  '
Set-TextValue "D49" '2.76'
Set-TextValue "E49" '  -1.09%  '
Set-TextValue "D50" '45.35'
Set-TextValue "E50" '  -5.20%  '
Set-TextValue "D51" '11.96'
Set-TextValue "E51" '  +18.96%  '
